$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is a plain decimal (e.g. "0.999", "599.50") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers and
# strips significant trailing/leading zeros (e.g. "599.50" -> 599.5).

$ws.Range("D2").Value = '66.810.86'
$ws.Range("E2").Value = '  -2.39%  '

$ws.Range("D3").Value = '3.471.32'
$ws.Range("E3").Value = '  -2.77%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.50'
$ws.Range("E5").Value = '  -3.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.81'
$ws.Range("E6").Value = '  -5.35%  '

$ws.Range("D7").Value = '3.473.42'
$ws.Range("E7").Value = '  -2.71%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  -2.79%  '

$ws.Range("E10").Value = '  -3.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.58'
$ws.Range("E11").Value = '  +3.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.421'
$ws.Range("E12").Value = '  -4.01%  '

$ws.Range("E13").Value = '  -4.50%  '

$ws.Range("D14").Value = '4.056.29'
$ws.Range("E14").Value = '  -2.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.26'
$ws.Range("E15").Value = '  -5.70%  '

$ws.Range("D16").Value = '3.472.70'
$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("D17").Value = '66.832.07'
$ws.Range("E17").Value = '  -2.25%  '

$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  -5.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.25'
$ws.Range("E20").Value = '  -4.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.97'
$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.14'
$ws.Range("E22").Value = '  -5.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.605'
$ws.Range("E23").Value = '  -6.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.91'
$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").Value = '3.607.82'
$ws.Range("E26").Value = '  -2.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -8.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.74'
$ws.Range("E28").Value = '  -7.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.34'
$ws.Range("E29").Value = '  -8.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.47'
$ws.Range("E30").Value = '  -3.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.57'
$ws.Range("E31").Value = '  -7.85%  '

$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("E33").Value = '  -3.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.20'
$ws.Range("E34").Value = '  -3.71%  '

$ws.Range("D35").Value = '3.460.21'
$ws.Range("E35").Value = '  -2.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.90'
$ws.Range("E36").Value = '  -6.98%  '

$ws.Range("E37").Value = '  -7.10%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.84'
$ws.Range("E39").Value = '  -5.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.87'
$ws.Range("E41").Value = '  -2.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0877'
$ws.Range("E42").Value = '  -4.66%  '

$ws.Range("E43").Value = '  -11.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.37'
$ws.Range("E44").Value = '  -4.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.891'
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.26'
$ws.Range("E46").Value = '  -0.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.55'
$ws.Range("E47").Value = '  -7.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.23'
$ws.Range("E48").Value = '  -8.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.41'
$ws.Range("E49").Value = '  -5.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.41'
$ws.Range("E50").Value = '  -9.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.971'
$ws.Range("E51").Value = '  -4.89%  '
